# Authoring.xlsx edit: "Y" -> "N" for the bulk of the "Runmode" column on the
# "Test Cases" sheet (row 2 is left as "Y"), normalize D76/D77 to the common
# cell style used by the rest of the column, move the active sheet / tab
# selection from "PostProfanityWordCheckTest" back to "Test Cases", and
# update the remembered scroll position / selection on both sheets.

$wb = $excel.ActiveWorkbook

$testCases = $wb.Worksheets.Item("Test Cases")
$postProfanity = $wb.Worksheets.Item("PostProfanityWordCheckTest")

# --- Test Cases!D3:D77 : "Y" -> "N" (D2 keeps its original "Y") ---
for ($r = 3; $r -le 77; $r++) {
    $testCases.Cells.Item($r, 4).Value = "N"
}

# D76 / D77 previously used slightly different (duplicate) cell styles;
# align them with the style used by the rest of the column (same style as D75).
$testCases.Range("D76").Style = $testCases.Range("D75").Style
$testCases.Range("D77").Style = $testCases.Range("D75").Style

# --- Restore "Test Cases" as the active / selected sheet & tab ---
$testCases.Activate()

# Scroll position + selection on "Test Cases" move one row further down.
$testCases.Application.ActiveWindow.ScrollRow = 49
$testCases.Range("D3:D77").Select()

# "PostProfanityWordCheckTest" is no longer the active tab; its remembered
# selection moves to D15.
$postProfanity.Range("D15").Select()
